$wb = $excel.ActiveWorkbook

# --- Update the "Hoja1" sheet text (A1) with the new conversion rates ---
$ws1 = $wb.Worksheets.Item("Hoja1")
$ws1.Range("A1").Value = "Conversión del día 💰`n✅ Dólar paralelo: 68`n`nBinance`n✅ 1000 Bs = 5.57 = 22512.51 pesos`n✅ 22512.51 pesos = 5.55 = 958.24 Bs`n`nPromedio competencia`n✅ Tasa pesos: 20`n✅ Tasa Bs: 20`n✅ % Ganancia: 20%"

# --- Update the "tasas" sheet numeric rate cells ---
$ws2 = $wb.Worksheets.Item("tasas")
$ws2.Range("N10").Value = 179.498
$ws2.Range("O10").Value = 4040.95
$ws2.Range("N12").Value = 4055
$ws2.Range("O12").Value = 172.6
